$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1: "Equip Log Missing Cost Code" data rows (4-12) ---
# Row 12 is a brand-new data row added at the end of this section.
$ws.Range("A4").Value  = "- 225010 Gilberto Ortiz 2026-02-13 1009"
$ws.Range("A5").Value  = "- 225010  2026-02-13 1009"
$ws.Range("A6").Value  = "- 225034  2026-02-12 1027"
$ws.Range("A7").Value  = "- 225034  2026-02-11 1027"
$ws.Range("A8").Value  = "- 225034  2026-02-10 1027"
$ws.Range("A9").Value  = "- 225034  2026-02-13 1027"
$ws.Range("A10").Value = "- 225034 Pablo Marin 2026-02-10 2014"
$ws.Range("A11").Value = "- 225034 Pablo Marin 2026-02-10 3031"
$ws.Range("A12").Value = "- 224034 Salvador Ortiz 2026-02-13 1013"

# --- The bold "section header" row for section 2 shifts from row 13 down to row 14 ---
# Copy the header formatting (bold font, same style as the other section headers) onto A14,
# then give A14 its header text, and blank out the old header row (13).
$ws.Range("A3").Copy()
$ws.Range("A14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A14").Value = "Equipment log entry with no matching time card entry"
$ws.Range("A13").Clear()

# --- Section 2: "Equipment log entry with no matching time card entry" data rows (15-22) ---
$ws.Range("A15").Value = "- 225010 Gilberto Ortiz 2026-02-11 200/500 1009"
$ws.Range("A16").Value = "- 225010 Gilberto Ortiz 2026-02-09 200/500 1009"
$ws.Range("A17").Value = "- 225010 Gilberto Ortiz 2026-02-10 200/500 1009"
$ws.Range("A18").Value = "- 225010 Doug Richards 2026-02-13 200/500 3026"
$ws.Range("A19").Value = "- 225010 Agustin Avila 2026-02-09 200/200 3026"
$ws.Range("A20").Value = "- 225010 Agustin Avila 2026-02-10 200/200 3042"
$ws.Range("A21").Value = "- 225010 Jesus Garcia 2026-02-13 200/500 1042"
$ws.Range("A22").Value = "- 225010 Salvador Martinez 2026-02-13 200/310 1010"

# --- The bold "section header" row for section 3 shifts from row 23 down to row 24 ---
$ws.Range("A3").Copy()
$ws.Range("A24").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A24").Value = "Equipment log with no operator"
$ws.Range("A23").Clear()

# --- Section 3: "Equipment log with no operator" data rows (25-29) ---
# Rows 26-29 are brand-new data rows added at the end of this section.
$ws.Range("A25").Value = "- Equip #: 225010 Job: 1009, Date: 2026-02-13"
$ws.Range("A26").Value = "- Equip #: 225034 Job: 1027, Date: 2026-02-12"
$ws.Range("A27").Value = "- Equip #: 225034 Job: 1027, Date: 2026-02-11"
$ws.Range("A28").Value = "- Equip #: 225034 Job: 1027, Date: 2026-02-10"
$ws.Range("A29").Value = "- Equip #: 225034 Job: 1027, Date: 2026-02-13"
